$wb = $excel.ActiveWorkbook

# --- ActorTable: insert new column "battltMusicOverriding|String" before the
#     last column (nodeWarLastCount|Int), which shifts from AB to AC. ---
$wsActor = $wb.Worksheets.Item("ActorTable")
$wsActor.Columns.Item(28).Insert()
$wsActor.Range("AB1").Value = "battltMusicOverriding|String"

# --- FixedCharTable: the group<->actorId rows 2-3 and 4-5 were swapped. ---
$wsFixed = $wb.Worksheets.Item("FixedCharTable")

$b2 = $wsFixed.Range("B2").Value()
$b3 = $wsFixed.Range("B3").Value()
$wsFixed.Range("B2").Value = $b3
$wsFixed.Range("B3").Value = $b2

$b4 = $wsFixed.Range("B4").Value()
$b5 = $wsFixed.Range("B5").Value()
$wsFixed.Range("B4").Value = $b5
$wsFixed.Range("B5").Value = $b4

# --- Active sheet moves from WingLookTable to ActorTable. ---
$wsActor.Activate()
$wsActor.Range("A1").Select()
